$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day text message ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.7 = 10002.97 pesos`n✅ 10002.97 pesos = 2.69 = 942.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the updated rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 370.4
$wsTasas.Range("O10").Value = 3705.1
$wsTasas.Range("O12").Value = 350.617
